# The workbook tracks, per column, a city name (row 1 header) plus one
# frequency value per year (rows 2-10). The edit re-groups a handful of
# city columns so that related/confusable city names end up adjacent,
# without touching any other data. For each group below, the columns'
# contents (header text + all data rows) are rotated right by one slot:
# the last column's values move into the first column's slot, and every
# other column shifts one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    ,@("B", "C")
    ,@("G", "H", "I")
    ,@("K", "L")
    ,@("U", "V")
    ,@("Y", "Z")
    ,@("AW", "AX")
    ,@("CL", "CM")
    ,@("DI", "DJ")
    ,@("DT", "DU")
)

$firstRow = 1
$lastRow = 10

foreach ($group in $groups) {
    $n = $group.Length

    for ($row = $firstRow; $row -le $lastRow; $row++) {
        # Snapshot the current values for every column in this group before
        # writing anything back, so later writes don't clobber values that
        # still need to be read.
        $vals = @()
        foreach ($col in $group) {
            $vals += , ($ws.Range($col + $row).Value())
        }

        for ($i = 0; $i -lt $n; $i++) {
            $srcIndex = ($i - 1 + $n) % $n
            $destCol = $group[$i]
            $ws.Range($destCol + $row).Value = $vals[$srcIndex]
        }
    }
}
